# The commit swaps the bodies of ppt/theme/theme1.xml (the slide-master
# theme, currently the "Simple Light" palette) and ppt/theme/theme2.xml
# (the notes-master theme, currently the unnamed "Default" palette) -
# theme1.xml ends up holding the old "Default" colors and theme2.xml ends
# up holding the old "Simple Light" colors. Font scheme / format scheme are
# identical between the two parts, so the only observable difference is the
# <a:clrScheme> (12 theme colors).
#
# Apply it through the Theme/ColorScheme object model: push the "Default"
# palette's RGB values onto the presentation's theme color scheme (and,
# for good measure, the color scheme reachable from the notes master too).

$p = $ppt.ActivePresentation

# RGB(...) in PowerPoint COM is 0x00BBGGRR, i.e. a little-endian BGR int.
function HexToRgbInt([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

# Target palette ("Default" scheme) in <a:clrScheme> child order.
$defaultScheme = [ordered]@{
    dk1      = "000000"
    lt1      = "FFFFFF"
    dk2      = "158158"
    lt2      = "F3F3F3"
    accent1  = "058DC7"
    accent2  = "50B432"
    accent3  = "ED561B"
    accent4  = "EDEF00"
    accent5  = "24CBE5"
    accent6  = "64E572"
    hlink    = "2200CC"
    folHlink = "551A8B"
}

function ApplySchemeToThemeColorScheme($themeColorScheme) {
    $i = 1
    foreach ($key in $defaultScheme.Keys) {
        $themeColorScheme.Item($i).RGB = HexToRgbInt $defaultScheme[$key]
        $i++
    }
}

# Primary path: the slide master's theme (ppt/theme/theme1.xml).
$master = $p.SlideMaster
ApplySchemeToThemeColorScheme $master.Theme.ThemeColorScheme

# Also push through the notes master's theme object for completeness -
# harmless if it resolves to the same underlying theme part.
if ($p.HasNotesMaster) {
    $notesMaster = $p.NotesMaster
    ApplySchemeToThemeColorScheme $notesMaster.Theme.ThemeColorScheme
}
